$d = $word.ActiveDocument

# 1. Replace the ID placeholder text (merges the two runs - text run + trailing
#    space run - into a single run with the updated placeholder text).
$d.Content.Find.Execute("**ID__AFFARS_5323_topic_3__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5323_370_4__ID**", 2)

# 2. Update the first paragraph's formatting: add a paragraph border (5pt
#    spacing on each side) and change the left indent from 6pt (120 twips)
#    to 11.25pt (225 twips).
$p = $d.Paragraphs(1)
$p.Borders.DistanceFromTop = 5
$p.Borders.DistanceFromLeft = 5
$p.Borders.DistanceFromBottom = 5
$p.Borders.DistanceFromRight = 5
$p.LeftIndent = 11.25
